$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we are about to update so that
# numeric-looking strings (e.g. "33.00", "0.530") keep their exact
# textual representation instead of being coerced to numbers.
$cellsToUpdate = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7",
    "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12",
    "E13", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18",
    "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24",
    "E24", "E25", "D26", "E26", "E27", "D28", "E28", "E29", "E30", "D31",
    "E31", "D32", "E32", "E33", "E34", "D35", "E35", "D36", "E36", "D37",
    "E37", "E38", "E39", "E40", "D41", "E41", "E42", "D43", "E43", "D44",
    "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49",
    "D50", "E50", "E51"
)
foreach ($cellRef in $cellsToUpdate) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.146.38'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '2.265.63'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '140.98'
$ws.Range("E5").Value = '  +13,984.89%  '
$ws.Range("D6").Value = '305.77'
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("D7").Value = '93.66'
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.486'
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("D11").Value = '33.00'
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("D12").Value = '0.0803'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").Value = '2.618.45'
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '14.36'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '2.269.21'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '0.784'
$ws.Range("E18").Value = '  +3.01%  '
$ws.Range("D19").Value = '41.990.37'
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = '12.68'
$ws.Range("E20").Value = '  +4.59%  '
$ws.Range("D21").Value = '0.0₃0916'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = '68.13'
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").Value = '243.67'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("E25").Value = '  +1.77%  '
$ws.Range("D26").Value = '1.94'
$ws.Range("E26").Value = '  +2.05%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '24.01'
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '35.04'
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("D32").Value = '159.68'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("E33").Value = '  +3.28%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").Value = '0.0744'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").Value = '3.06'
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").Value = '17.10'
$ws.Range("E37").Value = '  +2.51%  '
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("E42").Value = '  +3.77%  '
$ws.Range("D43").Value = '2.005.44'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("D44").Value = '19.48'
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("E45").Value = '  +10.16%  '
$ws.Range("D46").Value = '0.0284'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = '10.21'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '2.90'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = '53.60'
$ws.Range("E49").Value = '  +3.38%  '
$ws.Range("D50").Value = '73.10'
$ws.Range("E50").Value = '  +3.17%  '
$ws.Range("E51").Value = '  -0.24%  '

# Restore the default "Normal" style so no stray number-format
# style index is left behind on these cells.
foreach ($cellRef in $cellsToUpdate) {
    $ws.Range($cellRef).Style = "Normal"
}
